$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Some "Price" (column D) values are digit-grouped / decimal strings that
# Excel would otherwise auto-coerce to numbers (e.g. "1.00" -> 1, dropping
# the trailing zero). Mark those specific cells as Text before writing so
# the literal string is preserved, matching how the sheet already stores
# every Price cell as text.
$priceRows = @(2,3,5,6,10,12,13,14,15,16,17,20,21,22,23,24,25,26,27,28,29,30,31,32,34,35,39,40,41,42,43,44,46,47,48,49,50,51)
foreach ($r in $priceRows) {
    $ws.Range("D$r").NumberFormat = "@"
}

# Row 2
$ws.Range("D2").Value = "71.648.37"
$ws.Range("E2").Value = "  +4.44%  "

# Row 3
$ws.Range("D3").Value = "4.010.43"
$ws.Range("E3").Value = "  +4.37%  "

# Row 4
$ws.Range("E4").Value = "  +0.02%  "

# Row 5
$ws.Range("D5").Value = "536.44"
$ws.Range("E5").Value = "  +3.89%  "

# Row 6
$ws.Range("D6").Value = "152.04"
$ws.Range("E6").Value = "  +8.49%  "

# Row 7
$ws.Range("E7").Value = "  +14.22%  "

# Row 8
$ws.Range("E8").Value = "  +0.05%  "

# Row 9
$ws.Range("E9").Value = "  +5.81%  "

# Row 10
$ws.Range("D10").Value = "0.172"
$ws.Range("E10").Value = "  +2.71%  "

# Row 11
$ws.Range("E11").Value = "  +1.37%  "

# Row 12
$ws.Range("D12").Value = "47.95"
$ws.Range("E12").Value = "  +15.74%  "

# Row 13
$ws.Range("D13").Value = "10.73"
$ws.Range("E13").Value = "  +4.75%  "

# Row 14
$ws.Range("D14").Value = "4.647.83"
$ws.Range("E14").Value = "  +4.22%  "

# Row 15
$ws.Range("D15").Value = "4.013.83"
$ws.Range("E15").Value = "  +4.84%  "

# Row 16
$ws.Range("D16").Value = "14.08"
$ws.Range("E16").Value = "  +1.12%  "

# Row 17
$ws.Range("D17").Value = "20.54"
$ws.Range("E17").Value = "  -2.52%  "

# Row 18
$ws.Range("E18").Value = "  -0.33%  "

# Row 19
$ws.Range("E19").Value = "  +0.10%  "

# Row 20
$ws.Range("D20").Value = "71.532.78"
$ws.Range("E20").Value = "  +4.39%  "

# Row 21
$ws.Range("D21").Value = "430.76"
$ws.Range("E21").Value = "  +4.22%  "

# Row 22
$ws.Range("D22").Value = "98.52"
$ws.Range("E22").Value = "  +14.00%  "

# Row 23
$ws.Range("D23").Value = "3.53"
$ws.Range("E23").Value = "  +2.47%  "

# Row 24
$ws.Range("B24").Value = "PancakeSwap"
$ws.Range("C24").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D24").Value = "4.20"
$ws.Range("E24").Value = "  +5.27%  "

# Row 25
$ws.Range("B25").Value = "InternetComputer(DFINITY)"
$ws.Range("C25").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D25").Value = "14.48"
$ws.Range("E25").Value = "  +4.07%  "

# Row 26
$ws.Range("D26").Value = "11.11"
$ws.Range("E26").Value = "  -8.64%  "

# Row 27
$ws.Range("D27").Value = "10.73"
$ws.Range("E27").Value = "  +3.71%  "

# Row 28
$ws.Range("D28").Value = "3.72"
$ws.Range("E28").Value = "  +30.89%  "

# Row 29
$ws.Range("D29").Value = "5.84"
$ws.Range("E29").Value = "  +2.88%  "

# Row 30
$ws.Range("D30").Value = "36.78"
$ws.Range("E30").Value = "  +4.39%  "

# Row 31
$ws.Range("B31").Value = "Cosmos"
$ws.Range("C31").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D31").Value = "13.39"
$ws.Range("E31").Value = "  +0.42%  "

# Row 32
$ws.Range("B32").Value = "Bittensor"
$ws.Range("C32").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D32").Value = "689.38"
$ws.Range("E32").Value = "  +2.00%  "

# Row 33
$ws.Range("E33").Value = "  +5.19%  "

# Row 34
$ws.Range("D34").Value = "6.89"
$ws.Range("E34").Value = "  -1.02%  "

# Row 35
$ws.Range("D35").Value = "65.76"
$ws.Range("E35").Value = "  -0.72%  "

# Row 36
$ws.Range("E36").Value = "  +7.02%  "

# Row 37
$ws.Range("E37").Value = "  -4.03%  "

# Row 38
$ws.Range("E38").Value = "  +6.00%  "

# Row 39
$ws.Range("D39").Value = "3.53"
$ws.Range("E39").Value = "  +12.45%  "

# Row 40
$ws.Range("D40").Value = "0.0₃0826"
$ws.Range("E40").Value = "  -1.73%  "

# Row 41
$ws.Range("D41").Value = "3.46"
$ws.Range("E41").Value = "  +2.42%  "

# Row 42
$ws.Range("D42").Value = "1.00"
$ws.Range("E42").Value = "  +0.14%  "

# Row 43
$ws.Range("D43").Value = "0.998"
$ws.Range("E43").Value = "  -0.27%  "

# Row 44
$ws.Range("D44").Value = "0.0486"
$ws.Range("E44").Value = "  +3.04%  "

# Row 45
$ws.Range("E45").Value = "  +6.53%  "

# Row 46
$ws.Range("D46").Value = "2.63"
$ws.Range("E46").Value = "  -7.41%  "

# Row 47
$ws.Range("B47").Value = "THORChain"
$ws.Range("C47").Value = "https://coinranking.com/coin/ybmU-kKU+thorchain-rune"
$ws.Range("D47").Value = "9.52"
$ws.Range("E47").Value = "  +9.62%  "

# Row 48
$ws.Range("B48").Value = "ApeXProtocol"
$ws.Range("C48").Value = "https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex"
$ws.Range("D48").Value = "3.37"
$ws.Range("E48").Value = "  -2.36%  "

# Row 49
$ws.Range("D49").Value = "3.01"
$ws.Range("E49").Value = "  +0.86%  "

# Row 50
$ws.Range("D50").Value = "3.32"
$ws.Range("E50").Value = "  +1.20%  "

# Row 51
$ws.Range("D51").Value = "144.19"
$ws.Range("E51").Value = "  +0.99%  "
